# Apply button is now working properly, Motor Tuned with For Resp 30 and 25 and Tidal 700
#
# Adds a second worksheet ("Sheet2") after the existing "Sheet1", fills it
# with the tuning figures, formats the second value (row 2) in bold and the
# fourth value (row 4) with the scientific-notation number format already
# used elsewhere in the workbook, sets portrait page orientation, and makes
# Sheet2 the active/selected sheet (which also clears the old tab selection
# on Sheet1).

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)

# Insert the new sheet right after Sheet1 so tab order is Sheet1, Sheet2.
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "Sheet2"

# Tidal / resp tuning values.
$ws2.Range("A1").Value = 40
$ws2.Range("A2").Value = 30
$ws2.Range("A3").Value = 20
$ws2.Range("A4").Value = 0.02
$ws2.Range("A5").Value = 0.6
$ws2.Range("A6").Value = 0.8
$ws2.Range("A7").Value = 0.15

# Highlight the "Resp 30" figure in bold.
$ws2.Range("A2").Font.Bold = $true

# Match the scientific-notation format already used on Sheet1 for small values.
$ws2.Range("A4").NumberFormat = "0.00E+00"

# Portrait page for printing.
$ws2.PageSetup.Orientation = 1

# Make Sheet2 the active sheet/tab with B1 selected.
$ws2.Activate()
$ws2.Range("B1").Select()
